# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the values recorded at the newer data snapshot (456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for sheet "展览"
$zhanlanUpdates = @{
    2  = 232
    3  = 1077
    5  = 397
    6  = 73
    8  = 52
    9  = 6690
    10 = 142
    15 = 1069
    16 = 16009
    17 = 1573
    22 = 11252
    23 = 829
    24 = 4424
    25 = 291
    26 = 383
    27 = 37
    28 = 13
}

# Map of row -> new F value for sheet "全部类型"
$quanbuUpdates = @{
    2  = 232
    3  = 1077
    5  = 397
    6  = 73
    9  = 52
    10 = 6690
    11 = 142
    17 = 1069
    18 = 16009
    19 = 1573
    25 = 11252
    26 = 829
    27 = 4424
    28 = 291
    29 = 383
    30 = 37
    31 = 13
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Range("F$row").Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Range("F$row").Value = $quanbuUpdates[$row]
}
